# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Ultros_Profits workbook (profit calculation values per leve/sheet).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6255811
$ws.Range("I32").Value = 6872.875
$ws.Range("J32").Value = 12504749
$ws.Range("K32").Value = 6872.875
$ws.Range("L32").Value = 12504749
$ws.Range("M32").Value = -6546.875
$ws.Range("N32").Value = -12505401
$ws.Range("H39").Value = 203.11111
$ws.Range("I39").Value = 175.57143
$ws.Range("K39").Value = 526.71429
$ws.Range("M39").Value = -230.71429
$ws.Range("H58").Value = 2921.5715
$ws.Range("I58").Value = 290.4
$ws.Range("J58").Value = 9499.5
$ws.Range("K58").Value = 871.1999999999999
$ws.Range("L58").Value = 28498.5
$ws.Range("M58").Value = -721.1999999999999
$ws.Range("N58").Value = -28798.5
$ws.Range("H76").Value = 7666.3335
$ws.Range("H79").Value = 7666.3335
$ws.Range("H87").Value = 19999.953
$ws.Range("J87").Value = 19999.953
$ws.Range("L87").Value = 19999.953
$ws.Range("N87").Value = -22495.953
$ws.Range("H90").Value = 19999.953
$ws.Range("J90").Value = 19999.953
$ws.Range("L90").Value = 59999.859
$ws.Range("N90").Value = -72479.859
$ws.Range("H125").Value = 3291.4443
$ws.Range("I125").Value = 1430.75
$ws.Range("K125").Value = 12876.75
$ws.Range("M125").Value = -10416.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13169.053
$ws.Range("I32").Value = 9934.313
$ws.Range("K32").Value = 9934.313
$ws.Range("M32").Value = -9647.313
$ws.Range("H63").Value = 5275.5
$ws.Range("I63").Value = 2652
$ws.Range("J63").Value = 7899
$ws.Range("K63").Value = 2652
$ws.Range("L63").Value = 7899
$ws.Range("M63").Value = -1966
$ws.Range("N63").Value = -9271
$ws.Range("H66").Value = 5275.5
$ws.Range("I66").Value = 2652
$ws.Range("J66").Value = 7899
$ws.Range("K66").Value = 13260
$ws.Range("L66").Value = 39495
$ws.Range("M66").Value = -9828
$ws.Range("N66").Value = -46359
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41802
$ws.Range("H98").Value = 58632
$ws.Range("J98").Value = 58632
$ws.Range("L98").Value = 58632
$ws.Range("N98").Value = -64622
$ws.Range("H132").Value = 6229.0713
$ws.Range("I132").Value = 6232.316
$ws.Range("K132").Value = 18696.948
$ws.Range("M132").Value = -16166.948

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 13518.5
$ws.Range("I75").Value = 9316.5
$ws.Range("K75").Value = 9316.5
$ws.Range("M75").Value = -8380.5
$ws.Range("H78").Value = 13518.5
$ws.Range("I78").Value = 9316.5
$ws.Range("K78").Value = 27949.5
$ws.Range("M78").Value = -23269.5
$ws.Range("H99").Value = 37634.363
$ws.Range("I99").Value = 51101.75
$ws.Range("K99").Value = 51101.75
$ws.Range("M99").Value = -49603.75
$ws.Range("H100").Value = 38500.5
$ws.Range("J100").Value = 38500.5
$ws.Range("L100").Value = 38500.5
$ws.Range("N100").Value = -40664.5
$ws.Range("H107").Value = 2499.75
$ws.Range("I107").Value = 1999.5
$ws.Range("K107").Value = 1999.5
$ws.Range("M107").Value = -79.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 727.75
$ws.Range("I16").Value = 727.75
$ws.Range("K16").Value = 727.75
$ws.Range("M16").Value = -440.75
$ws.Range("H28").Value = 36216.715
$ws.Range("J28").Value = 36216.715
$ws.Range("L28").Value = 36216.715
$ws.Range("N28").Value = -36706.715
$ws.Range("H31").Value = 3856.625
$ws.Range("I31").Value = 2669.4827
$ws.Range("J31").Value = 6986.364
$ws.Range("K31").Value = 2669.4827
$ws.Range("L31").Value = 6986.364
$ws.Range("M31").Value = -2374.4827
$ws.Range("N31").Value = -7576.364
$ws.Range("H34").Value = 3856.625
$ws.Range("I34").Value = 2669.4827
$ws.Range("J34").Value = 6986.364
$ws.Range("K34").Value = 2669.4827
$ws.Range("L34").Value = 6986.364
$ws.Range("M34").Value = -2467.4827
$ws.Range("N34").Value = -7390.364
$ws.Range("H43").Value = 32106
$ws.Range("J43").Value = 32106
$ws.Range("L43").Value = 32106
$ws.Range("N43").Value = -32474
$ws.Range("H62").Value = 8997.857
$ws.Range("J62").Value = 8246.25
$ws.Range("L62").Value = 8246.25
$ws.Range("N62").Value = -9494.25
$ws.Range("H65").Value = 8997.857
$ws.Range("J65").Value = 8246.25
$ws.Range("L65").Value = 41231.25
$ws.Range("N65").Value = -47471.25
$ws.Range("H92").Value = 51525.25
$ws.Range("J92").Value = 51525.25
$ws.Range("L92").Value = 51525.25
$ws.Range("N92").Value = -56517.25
$ws.Range("H101").Value = 32106
$ws.Range("J101").Value = 32106
$ws.Range("L101").Value = 32106
$ws.Range("N101").Value = -38596
$ws.Range("H113").Value = 727.75
$ws.Range("I113").Value = 727.75
$ws.Range("K113").Value = 727.75
$ws.Range("M113").Value = 1442.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 41667780
$ws.Range("I22").Value = 832.6667
$ws.Range("J22").Value = 66667948
$ws.Range("K22").Value = 2498.0001
$ws.Range("L22").Value = 200003844
$ws.Range("M22").Value = -2329.0001
$ws.Range("N22").Value = -200004182
$ws.Range("H27").Value = 41667780
$ws.Range("I27").Value = 832.6667
$ws.Range("J27").Value = 66667948
$ws.Range("K27").Value = 2498.0001
$ws.Range("L27").Value = 200003844
$ws.Range("M27").Value = -2396.0001
$ws.Range("N27").Value = -200004048
$ws.Range("H29").Value = 27778628
$ws.Range("I29").Value = 1799.6
$ws.Range("K29").Value = 5398.799999999999
$ws.Range("M29").Value = -5121.799999999999
$ws.Range("H82").Value = 50002500
$ws.Range("I82").Value = 50002500
$ws.Range("K82").Value = 150007500
$ws.Range("M82").Value = -150007094
$ws.Range("H85").Value = 50002500
$ws.Range("I85").Value = 50002500
$ws.Range("K85").Value = 150007500
$ws.Range("M85").Value = -150006096
$ws.Range("H107").Value = 490.76666
$ws.Range("I107").Value = 444.2
$ws.Range("J107").Value = 500.08
$ws.Range("K107").Value = 1332.6
$ws.Range("L107").Value = 1500.24
$ws.Range("M107").Value = 587.4000000000001
$ws.Range("N107").Value = -5340.24
$ws.Range("H121").Value = 1830
$ws.Range("I121").Value = 261
$ws.Range("J121").Value = 2307.5217
$ws.Range("K121").Value = 783
$ws.Range("L121").Value = 6922.5651
$ws.Range("M121").Value = 527
$ws.Range("N121").Value = -9542.5651

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 226979.2
$ws.Range("I70").Value = 558555.5
$ws.Range("J70").Value = 5928.3335
$ws.Range("K70").Value = 558555.5
$ws.Range("L70").Value = 5928.3335
$ws.Range("M70").Value = -558285.5
$ws.Range("N70").Value = -6468.3335
$ws.Range("H73").Value = 226979.2
$ws.Range("I73").Value = 558555.5
$ws.Range("J73").Value = 5928.3335
$ws.Range("K73").Value = 558555.5
$ws.Range("L73").Value = 5928.3335
$ws.Range("M73").Value = -557619.5
$ws.Range("N73").Value = -7800.3335
$ws.Range("H86").Value = 20000
$ws.Range("J86").Value = 20000
$ws.Range("L86").Value = 20000
$ws.Range("N86").Value = -22372
$ws.Range("H89").Value = 20000
$ws.Range("J89").Value = 20000
$ws.Range("L89").Value = 60000
$ws.Range("N89").Value = -71856
$ws.Range("H98").Value = 16785.25
$ws.Range("J98").Value = 16785.25
$ws.Range("L98").Value = 16785.25
$ws.Range("N98").Value = -22775.25
$ws.Range("H107").Value = 283.4375
$ws.Range("I107").Value = 322.91666
$ws.Range("K107").Value = 322.91666
$ws.Range("M107").Value = 1597.08334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6495.8335
$ws.Range("I7").Value = 3325
$ws.Range("J7").Value = 7130
$ws.Range("K7").Value = 3325
$ws.Range("L7").Value = 7130
$ws.Range("M7").Value = -3213
$ws.Range("N7").Value = -7354
$ws.Range("H55").Value = 1704.6471
$ws.Range("I55").Value = 1995.1538
$ws.Range("K55").Value = 1995.1538
$ws.Range("M55").Value = -1822.1538
$ws.Range("H82").Value = 1998.875
$ws.Range("J82").Value = 2464.3333
$ws.Range("L82").Value = 2464.3333
$ws.Range("N82").Value = -3186.3333
$ws.Range("H85").Value = 1998.875
$ws.Range("J85").Value = 2464.3333
$ws.Range("L85").Value = 2464.3333
$ws.Range("N85").Value = -4960.3333
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 18000
$ws.Range("J97").Value = 18000
$ws.Range("L97").Value = 18000
$ws.Range("N97").Value = -19982
$ws.Range("H126").Value = 6495.8335
$ws.Range("I126").Value = 3325
$ws.Range("J126").Value = 7130
$ws.Range("K126").Value = 9975
$ws.Range("L126").Value = 21390
$ws.Range("M126").Value = -7505
$ws.Range("N126").Value = -26330

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17547056
$ws.Range("I81").Value = 2511.4285
$ws.Range("J81").Value = 66671784
$ws.Range("K81").Value = 5022.857
$ws.Range("L81").Value = 133343568
$ws.Range("M81").Value = -3961.857
$ws.Range("N81").Value = -133345690
$ws.Range("H84").Value = 17547056
$ws.Range("I84").Value = 2511.4285
$ws.Range("J84").Value = 66671784
$ws.Range("K84").Value = 25114.285
$ws.Range("L84").Value = 666717840
$ws.Range("M84").Value = -19810.285
$ws.Range("N84").Value = -666728448
$ws.Range("H107").Value = 7730.6665
$ws.Range("I107").Value = 984.1429000000001
$ws.Range("K107").Value = 2952.4287
$ws.Range("M107").Value = -1032.4287
